$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The sheet currently lists 5 missing/low-stock items in rows 7-11, a totals
# row (12) and a footer row (13). The new version lists 8 items (3 more were
# added to the underlying report, re-sorted alphabetically), so 3 extra rows
# are needed before the totals/footer rows, which get pushed from 12/13 to
# 15/16.
# ---------------------------------------------------------------------------

# Insert 3 blank rows right above the old totals row (row 12). Excel shifts
# the totals row (12 -> 15), the footer row (13 -> 16), their merged ranges
# and row heights down automatically.
$ws.Rows("12:14").Insert()

# Clone the formatting (styles + column merges) of an existing data row
# (row 7) into each of the 3 newly-inserted rows so they look like the rest
# of the table.
for ($r = 12; $r -le 14; $r++) {
    for ($c = 1; $c -le 17; $c++) {
        $ws.Cells.Item(7, $c).Copy($ws.Cells.Item($r, $c))
    }
}
$excel.CutCopyMode = $false

foreach ($r in 12,13,14) {
    $ws.Range("A" + $r + ":B" + $r).Merge()
    $ws.Range("C" + $r + ":G" + $r).Merge()
    $ws.Range("H" + $r + ":K" + $r).Merge()
    $ws.Range("L" + $r + ":M" + $r).Merge()
    $ws.Range("N" + $r + ":O" + $r).Merge()
}

# Row heights follow the alternating 25.5 / 24.75 pattern already used by
# rows 7-11, and the relocated totals row (now 15) becomes 24.75 too.
$ws.Rows("12").RowHeight = 25.5
$ws.Rows("13").RowHeight = 24.75
$ws.Rows("14").RowHeight = 25.5
$ws.Rows("15").RowHeight = 24.75

# ---------------------------------------------------------------------------
# Write the final (re-sorted) list of 8 items into rows 7-14.
# Columns: A = row #, C = item name, H = current/reorder qty, L = count,
#          N = price, P = sale price, Q = transaction count.
# ---------------------------------------------------------------------------
$items = @(
    @{ A = 1; C = "BETADERM 0.1% CREAM 15 GM";          H = "4:0";    L = "1"; N = "18.00"; P = "18.0000"; Q = "1:0" },
    @{ A = 2; C = "CETAL 250MG/5ML 60ML SUSP";           H = "8:0";    L = "1"; N = "31.00"; P = "31.0000"; Q = "1:0" },
    @{ A = 3; C = "CONGESTAL 20 TABS";                   H = "0:1";    L = "1"; N = "50.00"; P = "25.0000"; Q = "0:1" },
    @{ A = 4; C = "FLAGYL 125MG/5ML 100 ML SUSPENSION";  H = "7:0";    L = "1"; N = "26.00"; P = "26.0000"; Q = "1:0" },
    @{ A = 5; C = "GARAMYCIN 0.1% OINT. 15 GM";          H = "0:0";    L = "1"; N = "22.00"; P = "22.0000"; Q = "1:0" },
    @{ A = 6; C = "HIBIOTIC N 457MG/5ML SUSP. 60ML";     H = "1:0";    L = "1"; N = "80.00"; P = "80.0000"; Q = "1:0" },
    @{ A = 7; C = "HIDERM TOPICAL LOTION 100ML";         H = "0:0";    L = "1"; N = "55.00"; P = "55.0000"; Q = "1:0" },
    @{ A = 8; C = "WATER FOR INJECTION AMP. 5 ML";       H = "8403:0"; L = "1"; N = "2.00";  P = "4.0000";  Q = "2:0" }
)

$row = 7
foreach ($item in $items) {
    $ws.Range("A" + $row).Value = $item.A
    $ws.Range("C" + $row).Value = $item.C
    $ws.Range("H" + $row).Value = $item.H
    $ws.Range("L" + $row).Value = $item.L
    $ws.Range("N" + $row).Value = $item.N
    $ws.Range("P" + $row).Value = $item.P
    $ws.Range("Q" + $row).Value = $item.Q
    $row = $row + 1
}

# ---------------------------------------------------------------------------
# Update the grand-total cell (sum of the price column) and the generated
# timestamp in the footer.
# ---------------------------------------------------------------------------
$ws.Range("P15").Value = 261
$ws.Range("A16").Value = "Thursday, 14 August, 2025 10:01 AM"
